# Update the "Förändrad" (Changed) date column (C) for rows 2-67
# from serial date 45189 (2023-09-20) to serial date 45190 (2023-09-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 67; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
